$d = $word.ActiveDocument

# 1) Replace the "130" estado cell text with "Pendiente"
$d.Content.Find.Execute("130", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pendiente", 2)

# 2) Replace the "145" estado cell text with "Pendiente"
$d.Content.Find.Execute("145", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pendiente", 2)

# 3) Replace the task text for 0.4.2 row
$d.Content.Find.Execute("Probar predicción sobre CSV actual", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Integrar modelo en entorno de prueba/programa final.", 2)

# 4) Set the row height for the "0.4.1" row (row index 4, 1-based 5th row in table 1)
$table = $d.Tables.Item(1)
$row = $table.Rows.Item(5)
$row.HeightRule = 1
$row.Height = 197.37304687500006
